$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as genuine text (no numeric/date coercion, no style churn).
# Builds a temp formula returning the literal string, then Copy + PasteSpecial(values)
# collapses it back down to a plain shared-string cell, matching how the sheet was
# authored (t="inlineStr"/shared-string cells, no numFmt changes).
function Set-TextValue($addr, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "57.095.48"
Set-TextValue "E2" "  +1.77%  "

Set-TextValue "D3" "3.252.95"
Set-TextValue "E3" "  +0.97%  "

Set-TextValue "E4" "  +0.05%  "

Set-TextValue "D5" "398.18"
Set-TextValue "E5" "  -0.21%  "

Set-TextValue "D6" "108.66"
Set-TextValue "E6" "  -1.51%  "

Set-TextValue "D7" "0.579"
Set-TextValue "E7" "  +4.35%  "

Set-TextValue "E8" "  +0.00%  "

Set-TextValue "E9" "  -0.82%  "

Set-TextValue "D10" "39.27"
Set-TextValue "E10" "  -0.29%  "

Set-TextValue "E11" "  +5.12%  "

Set-TextValue "E12" "  +1.59%  "

Set-TextValue "D13" "3.772.46"
Set-TextValue "E13" "  +1.23%  "

Set-TextValue "D14" "8.25"
Set-TextValue "E14" "  +1.31%  "

Set-TextValue "D15" "18.96"
Set-TextValue "E15" "  -0.71%  "

Set-TextValue "D16" "3.250.90"
Set-TextValue "E16" "  +0.95%  "

Set-TextValue "E17" "  -2.61%  "

Set-TextValue "E18" "  +2.82%  "

Set-TextValue "D19" "56.911.57"
Set-TextValue "E19" "  +1.72%  "

Set-TextValue "D20" "3.30"
Set-TextValue "E20" "  -1.59%  "

Set-TextValue "E21" "  +5.10%  "

Set-TextValue "D22" "12.93"
Set-TextValue "E22" "  -1.43%  "

Set-TextValue "D23" "293.55"
Set-TextValue "E23" "  -3.95%  "

Set-TextValue "D24" "74.01"
Set-TextValue "E24" "  -1.74%  "

Set-TextValue "E25" "  -1.79%  "

Set-TextValue "E28" "  +0.64%  "

Set-TextValue "E29" "  -0.41%  "

Set-TextValue "E30" "  -2.51%  "

Set-TextValue "E31" "  -0.03%  "

Set-TextValue "E32" "  +0.79%  "

Set-TextValue "D33" "11.21"
Set-TextValue "E33" "  -0.52%  "

Set-TextValue "D34" "40.12"
Set-TextValue "E34" "  +10.78%  "

Set-TextValue "E35" "  -0.35%  "

Set-TextValue "E36" "  +0.88%  "

Set-TextValue "D37" "51.25"
Set-TextValue "E37" "  -0.26%  "

Set-TextValue "D38" "0.999"
Set-TextValue "E38" "  -0.03%  "

Set-TextValue "D39" "3.47"
Set-TextValue "E39" "  -0.93%  "

Set-TextValue "D40" "2.98"
Set-TextValue "E40" "  -3.98%  "

Set-TextValue "D41" "136.53"
Set-TextValue "E41" "  +1.25%  "

Set-TextValue "E42" "  +1.77%  "

Set-TextValue "E43" "  -0.83%  "

Set-TextValue "E44" "  -2.91%  "

Set-TextValue "D45" "3.91"
Set-TextValue "E45" "  -3.88%  "

Set-TextValue "D46" "16.80"
Set-TextValue "E46" "  -2.88%  "

Set-TextValue "E47" "  -0.27%  "

Set-TextValue "E48" "  +4.72%  "

Set-TextValue "D49" "2.141.86"
Set-TextValue "E49" "  +0.00%  "

Set-TextValue "D50" "2.45"
Set-TextValue "E50" "  -2.15%  "

Set-TextValue "D51" "1.96"
Set-TextValue "E51" "  -6.50%  "

# Rows 26 and 27: Filecoin and EthereumClassic swapped rank position, with updated values
Set-TextValue "B26" "EthereumClassic"
Set-TextValue "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D26" "28.02"
Set-TextValue "E26" "  -0.91%  "

Set-TextValue "B27" "Filecoin"
Set-TextValue "C27" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D27" "7.90"
Set-TextValue "E27" "  -4.23%  "

$excel.CutCopyMode = $false
